$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be read/written as text so values like
# "30.331.96" or "0.9994" are not auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.331.96'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '1.936.46'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = '251.43'
$ws.Range("E5").Value = '  +1.81%  '

$ws.Range("D6").Value = '0.7267'
$ws.Range("E6").Value = '  +3.39%  '

$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").Value = '  -0.53%  '

$ws.Range("D8").Value = '0.3316'
$ws.Range("E8").Value = '  +2.34%  '

$ws.Range("D9").Value = '28.05'
$ws.Range("E9").Value = '  +6.45%  '

$ws.Range("D10").Value = '0.07289'
$ws.Range("E10").Value = '  +6.90%  '

$ws.Range("D11").Value = '0.8104'
$ws.Range("E11").Value = '  +1.84%  '

$ws.Range("D12").Value = '0.08106'
$ws.Range("E12").Value = '  +1.81%  '

$ws.Range("D13").Value = '1.934.39'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").Value = '5.492'
$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("D15").Value = '94.99'
$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").Value = '15.16'
$ws.Range("E16").Value = '  +5.21%  '

$ws.Range("D17").Value = '30.324.25'
$ws.Range("E17").Value = '  +0.15%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.000008268'
$ws.Range("E18").Value = '  +5.97%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '254.27'
$ws.Range("E19").Value = '  -2.44%  '

$ws.Range("D20").Value = '5.827'
$ws.Range("E20").Value = '  -0.40%  '

$ws.Range("D21").Value = '2.187.85'
$ws.Range("E21").Value = '  +0.11%  '

$ws.Range("D22").Value = '0.9992'
$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").Value = '0.9991'
$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("D24").Value = '6.983'
$ws.Range("E24").Value = '  +2.63%  '

$ws.Range("D25").Value = '9.784'
$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("D26").Value = '165.51'
$ws.Range("E26").Value = '  +4.20%  '

$ws.Range("D27").Value = '2.359'
$ws.Range("E27").Value = '  +6.13%  '

$ws.Range("D28").Value = '19.38'
$ws.Range("E28").Value = '  +3.13%  '

$ws.Range("E29").Value = '  -0.45%  '

$ws.Range("E30").Value = '  +0.91%  '

$ws.Range("D31").Value = '1.538'
$ws.Range("E31").Value = '  -1.10%  '

$ws.Range("D32").Value = '4.448'
$ws.Range("E32").Value = '  +0.99%  '

$ws.Range("D33").Value = '4.212'
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("D34").Value = '0.05270'
$ws.Range("E34").Value = '  +4.03%  '

$ws.Range("D35").Value = '1.275'
$ws.Range("E35").Value = '  +6.91%  '

$ws.Range("D36").Value = '0.7511'
$ws.Range("E36").Value = '  +0.50%  '

$ws.Range("D37").Value = '2.764'
$ws.Range("E37").Value = '  +1.88%  '

$ws.Range("D38").Value = '0.01979'
$ws.Range("E38").Value = '  +3.10%  '

$ws.Range("D39").Value = '2.796'
$ws.Range("E39").Value = '  +0.96%  '

$ws.Range("D40").Value = '79.39'
$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("D41").Value = '6.443'
$ws.Range("E41").Value = '  -1.38%  '

$ws.Range("D42").Value = '0.4571'
$ws.Range("E42").Value = '  +3.60%  '

$ws.Range("D43").Value = '2.040'
$ws.Range("E43").Value = '  +0.31%  '

$ws.Range("D44").Value = '0.8455'
$ws.Range("E44").Value = '  +1.11%  '

$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.46%  '

$ws.Range("D46").Value = '101.95'
$ws.Range("E46").Value = '  +0.31%  '

$ws.Range("D47").Value = '9.790'
$ws.Range("E47").Value = '  +2.11%  '

$ws.Range("D48").Value = '7.478'
$ws.Range("E48").Value = '  +3.74%  '

$ws.Range("D49").Value = '36.72'
$ws.Range("E49").Value = '  +2.83%  '

$ws.Range("D50").Value = '0.4214'
$ws.Range("E50").Value = '  +3.87%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.504'
$ws.Range("E51").Value = '  +1.76%  '

# Restore the default (unformatted) style on the Price column now that
# the literal text has been written, matching the original workbook look.
$ws.Range("D2:D51").Style = "Normal"
